# Remove the (empty) footnote reference from the single paragraph and its
# corresponding footnote entry in the footnotes part. This fixes the
# "bug with empty notes" — the footnote body only ever contained the
# footnote mark plus a lone "།" with no real annotation text, so the
# whole footnote (reference + note) is deleted.

$d = $word.ActiveDocument

$footnotes = $d.Footnotes
$count = $footnotes.Count

for ($i = $count; $i -ge 1; $i--) {
    $footnotes.Item($i).Delete()
}

Write-Output "Remaining footnotes: $($d.Footnotes.Count)"
